$wb = $excel.ActiveWorkbook
$wsAluno = $wb.Worksheets.Item("Aluno")

# Mirror the "Professor" sheet's IFERROR formula into column H of the
# "Aluno" sheet, but falling back to 0 instead of the text "Atenção".
# Enter H6 first (own formula), then fill H7:H12 so Excel records it the
# same way as on the "Professor" sheet: H6 standalone, H7:H12 shared.
$wsAluno.Range("H6").Formula = "=IFERROR(F6+G6,0)"
$wsAluno.Range("H7:H12").Formula = "=IFERROR(F7+G7,0)"

# Make "Aluno" the active sheet / tab, with J7 selected.
$wsAluno.Activate()
$wsAluno.Range("J7").Select()
